$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 22 and 23 swap content (BitcoinCash <-> PancakeSwap) with updated values
$ws.Cells.Item(22, 2).Value = "PancakeSwap"
$ws.Cells.Item(22, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(22, 4).Formula = "'3.13"
$ws.Cells.Item(22, 5).Value = "  +5.72%  "

$ws.Cells.Item(23, 2).Value = "BitcoinCash"
$ws.Cells.Item(23, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(23, 4).Formula = "'237.39"
$ws.Cells.Item(23, 5).Value = "  -1.46%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
$ws.Cells.Item(2, 4).Value = "44.364.08"
$ws.Cells.Item(2, 5).Value = "  +0.44%  "
$ws.Cells.Item(3, 4).Value = "2.242.98"
$ws.Cells.Item(3, 5).Value = "  -0.66%  "
$ws.Cells.Item(4, 5).Value = "  +0.58%  "
$ws.Cells.Item(5, 4).Formula = "'306.26"
$ws.Cells.Item(5, 5).Value = "  -1.02%  "
$ws.Cells.Item(6, 4).Formula = "'93.26"
$ws.Cells.Item(6, 5).Value = "  -5.94%  "
$ws.Cells.Item(7, 4).Formula = "'0.572"
$ws.Cells.Item(7, 5).Value = "  -0.68%  "
$ws.Cells.Item(8, 5).Value = "  +0.31%  "
$ws.Cells.Item(9, 4).Formula = "'0.521"
$ws.Cells.Item(9, 5).Value = "  -3.05%  "
$ws.Cells.Item(10, 4).Formula = "'34.55"
$ws.Cells.Item(10, 5).Value = "  -3.41%  "
$ws.Cells.Item(11, 4).Formula = "'0.0809"
$ws.Cells.Item(11, 5).Value = "  -1.78%  "
$ws.Cells.Item(12, 4).Formula = "'7.14"
$ws.Cells.Item(12, 5).Value = "  -3.11%  "
$ws.Cells.Item(13, 5).Value = "  -0.05%  "
$ws.Cells.Item(14, 4).Value = "2.360.21"
$ws.Cells.Item(14, 5).Value = "  +4.47%  "
$ws.Cells.Item(15, 4).Formula = "'0.834"
$ws.Cells.Item(15, 5).Value = "  -1.13%  "
$ws.Cells.Item(16, 4).Formula = "'13.54"
$ws.Cells.Item(16, 5).Value = "  -2.81%  "
$ws.Cells.Item(17, 4).Value = "44.046.18"
$ws.Cells.Item(17, 5).Value = "  -0.06%  "
$ws.Cells.Item(18, 4).Value = "0.0₃0960"
$ws.Cells.Item(18, 5).Value = "  -1.72%  "
$ws.Cells.Item(19, 4).Formula = "'12.33"
$ws.Cells.Item(19, 5).Value = "  -4.79%  "
$ws.Cells.Item(20, 4).Formula = "'6.34"
$ws.Cells.Item(20, 5).Value = "  -0.51%  "
$ws.Cells.Item(21, 4).Formula = "'65.65"
$ws.Cells.Item(21, 5).Value = "  +0.10%  "
$ws.Cells.Item(24, 4).Formula = "'1.98"
$ws.Cells.Item(24, 5).Value = "  -0.60%  "
$ws.Cells.Item(25, 5).Value = "  -0.56%  "
$ws.Cells.Item(26, 4).Formula = "'38.74"
$ws.Cells.Item(26, 5).Value = "  +4.22%  "
$ws.Cells.Item(27, 5).Value = "  +3.94%  "
$ws.Cells.Item(28, 4).Formula = "'9.79"
$ws.Cells.Item(28, 5).Value = "  -4.09%  "
$ws.Cells.Item(29, 4).Formula = "'5.94"
$ws.Cells.Item(29, 5).Value = "  -4.41%  "
$ws.Cells.Item(30, 4).Formula = "'19.98"
$ws.Cells.Item(30, 5).Value = "  -0.98%  "
$ws.Cells.Item(31, 4).Formula = "'153.23"
$ws.Cells.Item(31, 5).Value = "  -2.98%  "
$ws.Cells.Item(32, 4).Formula = "'0.0796"
$ws.Cells.Item(32, 5).Value = "  -4.05%  "
$ws.Cells.Item(33, 5).Value = "  -0.96%  "
$ws.Cells.Item(34, 4).Formula = "'3.12"
$ws.Cells.Item(34, 5).Value = "  -12.63%  "
$ws.Cells.Item(35, 4).Formula = "'0.109"
$ws.Cells.Item(35, 5).Value = "  +1.16%  "
$ws.Cells.Item(36, 4).Formula = "'0.119"
$ws.Cells.Item(36, 5).Value = "  -0.44%  "
$ws.Cells.Item(37, 4).Formula = "'1.77"
$ws.Cells.Item(37, 5).Value = "  -4.96%  "
$ws.Cells.Item(38, 4).Formula = "'3.46"
$ws.Cells.Item(38, 5).Value = "  +1.34%  "
$ws.Cells.Item(39, 4).Formula = "'14.67"
$ws.Cells.Item(39, 5).Value = "  -7.76%  "
$ws.Cells.Item(40, 4).Formula = "'3.81"
$ws.Cells.Item(40, 5).Value = "  -2.81%  "
$ws.Cells.Item(41, 4).Formula = "'0.0301"
$ws.Cells.Item(41, 5).Value = "  -2.06%  "
$ws.Cells.Item(42, 5).Value = "  +0.48%  "
$ws.Cells.Item(43, 4).Value = "1.740.14"
$ws.Cells.Item(43, 5).Value = "  -2.03%  "
$ws.Cells.Item(44, 4).Formula = "'80.77"
$ws.Cells.Item(44, 5).Value = "  -7.59%  "
$ws.Cells.Item(45, 4).Formula = "'0.191"
$ws.Cells.Item(45, 5).Value = "  -1.88%  "
$ws.Cells.Item(46, 5).Value = "  +4.54%  "
$ws.Cells.Item(47, 4).Formula = "'99.44"
$ws.Cells.Item(47, 5).Value = "  -2.51%  "
$ws.Cells.Item(48, 4).Formula = "'4.92"
$ws.Cells.Item(48, 5).Value = "  -4.87%  "
$ws.Cells.Item(49, 4).Formula = "'14.72"
$ws.Cells.Item(49, 5).Value = "  +4.52%  "
$ws.Cells.Item(50, 4).Formula = "'8.16"
$ws.Cells.Item(50, 5).Value = "  -1.89%  "
$ws.Cells.Item(51, 4).Formula = "'55.16"
$ws.Cells.Item(51, 5).Value = "  -1.05%  "
